$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": update Segundo Parcial statistics row ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 1
$ws2P.Range("F2").Value = 16
$ws2P.Range("G2").Value = 94.12
$ws2P.Range("H2").Value = 8.4

# --- Sheet "Rescatables": update Reprobadas count for the student ---
$wsResc = $wb.Worksheets.Item("Rescatables")
$wsResc.Range("G2").Value = 3
